$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2-416). The commit updates that date from 2023-09-06 (serial 45175)
# to 2023-09-08 (serial 45177) for every row.
$oldSerial = 45175
$newSerial = 45177

$lastRow = 416
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq $oldSerial) {
        $cell.Value2 = $newSerial
    }
}
